$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated matching results (Logit-scale matching instead of probability scale)
$values = @{
    2 = @(6.687148029338182, 6.570311120900312, 0.11683690843787, 12.70713733727643, 0.9030117737368897, 0.3676060922413187)
    3 = @(5.848378376315575, 5.773717978679917, 0.07466039763565835, 8.525852126170705, 0.6058756308934319, 0.5452912898077318)
    4 = @(0.7197834411214429, 0.7471074021414038, -0.02732396101996093, -6.683957501942611, -0.474984424832292, 0.6353208663150816)
    5 = @(0.4945814727379544, 0.4770474889431849, 0.01753398379476945, 17.52308519848812, 1.245249171897352, 0.2145669291153803)
}

foreach ($row in $values.Keys) {
    $cols = $values[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = 2 + $i  # column B is index 2
        $ws.Cells.Item($row, $col).Value = $cols[$i]
    }
}

$wb.Save()
